$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the "arts" column (currently R) to make
# room for history / electives / cs subject columns.
$ws.Range("R1:T1").EntireColumn.Insert()

# New header cells for the inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New row-2 values for the inserted columns (numeric 0, like the other
# general_college_subjects.* columns).
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Lower-case the descriptive text values in row 2 for columns D..J.
$ws.Range("D2").Value = "not considered"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "not considered"
$ws.Range("G2").Value = "important"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"
